$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "20.264.80"
$ws.Range("D3").Value = "1.440.01"
$ws.Range("E3").Value = "  +2.31%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.005"
$ws.Range("D4").Style = "Normal"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.9128"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -8.80%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "274.97"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.54%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.3628"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -1.65%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3077"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -1.82%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "39.06"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -1.56%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.021"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +0.79%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.06504"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.12%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.9983"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +0.01%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "5.347"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -2.33%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "17.54"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +0.77%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "6.042"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -1.79%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.00001011"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -0.47%  "
$ws.Range("D17").Value = "1.438.28"
$ws.Range("E17").Value = "  +2.38%  "
$ws.Range("E18").Value = "  -6.96%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.05632"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -1.27%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "67.49"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -5.02%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "5.399"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -3.24%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "14.24"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -3.48%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "10.82"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -2.16%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.240"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.72%  "
$ws.Range("D25").Value = "20.263.88"
$ws.Range("E25").Value = "  +1.21%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "137.65"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +1.35%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.098"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -6.44%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "16.91"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.10%  "
$ws.Range("D29").Value = "1.592.03"
$ws.Range("E29").Value = "  +1.82%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "110.39"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +0.77%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "3.934"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -3.97%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.8046"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -2.12%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.850"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -8.77%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.07654"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -0.20%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.462"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +1.22%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.05816"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +0.52%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "4.677"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -3.32%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.131"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +2.37%  "
$ws.Range("B39").Value = "Aptos"
$ws.Range("C39").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "10.20"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -2.20%  "
$ws.Range("B40").Value = "VeChain"
$ws.Range("C40").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.01984"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -4.63%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.1852"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -2.66%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.9258"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -7.40%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "7.102"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -15.92%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.5211"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -1.47%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "3.485"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.78%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "11.86"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -3.27%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "116.69"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +4.40%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.5094"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.77%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.730"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -2.27%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.06417"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +3.91%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.9752"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -2.33%  "
